$d = $word.ActiveDocument
$d.Content.Find.Execute(" and animation. All topics were covered in Java", $true, $false, $false, $false, $false, $true, 1, $false, " and 2D animation. All topics were covered in Java", 2)
